# Fruta / hortaliza, semanal
# Insert a new weekly record at row 47 (shifting existing rows 47-108 down to 48-109)
# for "Hortaliza, Feria Lagunitas de Puerto Montt - Haba".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(47).Insert()

$ws.Cells.Item(47, 1).Value = 4
$ws.Cells.Item(47, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(47, 3).Value = "Los Lagos"
$ws.Cells.Item(47, 4).Value = 44848
$ws.Cells.Item(47, 5).Value = 10
$ws.Cells.Item(47, 6).Value = 100112026
$ws.Cells.Item(47, 7).Value = "Haba"
$ws.Cells.Item(47, 8).Value = "Sin especificar"
$ws.Cells.Item(47, 9).Value = "Primera"
$ws.Cells.Item(47, 10).Value = 180
$ws.Cells.Item(47, 11).Value = 10000
$ws.Cells.Item(47, 12).Value = 10000
$ws.Cells.Item(47, 13).Value = 10000
$ws.Cells.Item(47, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(47, 15).Value = "Región Metropolitana"
$ws.Cells.Item(47, 16).Value = 400
$ws.Cells.Item(47, 17).Value = 25
$ws.Cells.Item(47, 18).Value = "Hortaliza"
